$wb = $excel.ActiveWorkbook
$progress = $wb.Worksheets.Item(1)
$journal  = $wb.Worksheets.Item(2)

# --- Progress sheet -------------------------------------------------------
# Row 9: "Low-Pass & High-Pass Filtering" - end date moved earlier
$progress.Range("G9").Value = 45730

# Row 10: "Artifact Detection & Removal" - start date filled in
$progress.Range("F10").Value = 45730

# Move the visible selection on the Progress sheet to I10
$progress.Activate()
$progress.Range("I10").Select()

# --- Journal sheet ----------------------------------------------------------
# Update the long journal note about the high-pass filter conv() fix
$oldNote = $journal.Range("C5").Value()
$newNote = $oldNote.Replace(
    "- the conv function just did a 1d convolution with the signal samples and coeffs." + [char]10,
    "- the conv function just did a 1d convolution with the signal samples and coeffs. The conv function adds 50 samples to the signal length, which seems to be how convolution works." + [char]10
)
$newNote = $newNote.Replace(
    "recreate." + [char]10 + [char]10 + "Although",
    "recreate." + [char]10 + "Although"
)
$newNote = $newNote.Replace(
    "I wonder if this is why some misdetections and ",
    "I wonder if this is why some misdetections and glitches"
)
$journal.Range("C5").Value = $newNote

# New journal entry logging the fix
$journal.Range("A6").Value = "Artifact Detection"
$journal.Range("B6").Value = 45733

$journal.Activate()
$journal.Range("D5").Select()
